# Auto-generated Excel COM-interop script applying the weekly CompStat data refresh
# (new crime data collected) per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# A8: 'Volume 29   Number  50' -> 'Volume 29   Number  51'
$cellA8 = $ws.Range("A8")
$cellA8.Characters(21,2).Text = "51"

# C9: 'Report Covering the Week  12/12/2022  Through  12/18/2022'
#  -> 'Report Covering the Week  12/19/2022  Through  12/25/2022'
$cellC9 = $ws.Range("C9")
$cellC9.Characters(27,10).Text = "12/19/2022"
$cellC9.Characters(48,10).Text = "12/25/2022"

# --- Precinct / crime statistics table updates (rows 14-30) ---

# Row 14
$ws.Range("D14").Value2 = 1
$ws.Range("J14").Value2 = 93
$ws.Range("K14").Value2 = -18.279569892473
$ws.Range("L14").Value2 = -26.923076923076
$ws.Range("M14").Value2 = -41.538461538461
$ws.Range("N14").Value2 = -83.829787234042

# Row 15
$ws.Range("C15").Value2 = 9
$ws.Range("D15").Value2 = 2
$ws.Range("E15").Value2 = 350
$ws.Range("F15").Value2 = 19
$ws.Range("G15").Value2 = 15
$ws.Range("H15").Value2 = 26.666666666666
$ws.Range("I15").Value2 = 250
$ws.Range("J15").Value2 = 218
$ws.Range("K15").Value2 = 14.678899082568
$ws.Range("L15").Value2 = 13.122171945701
$ws.Range("M15").Value2 = 11.111111111111
$ws.Range("N15").Value2 = -57.555178268251

# Row 16
$ws.Range("C16").Value2 = 34
$ws.Range("D16").Value2 = 47
$ws.Range("E16").Value2 = -27.659574468085
$ws.Range("F16").Value2 = 157
$ws.Range("G16").Value2 = 188
$ws.Range("H16").Value2 = -16.489361702127
$ws.Range("I16").Value2 = 2495
$ws.Range("J16").Value2 = 2080
$ws.Range("K16").Value2 = 19.951923076923
$ws.Range("L16").Value2 = 17.744218971212
$ws.Range("M16").Value2 = -32.016348773842
$ws.Range("N16").Value2 = -84.987063000180

# Row 17
$ws.Range("C17").Value2 = 65
$ws.Range("D17").Value2 = 63
$ws.Range("E17").Value2 = 3.174603174603
$ws.Range("G17").Value2 = 280
$ws.Range("H17").Value2 = -9.285714285714
$ws.Range("I17").Value2 = 4029
$ws.Range("J17").Value2 = 3558
$ws.Range("K17").Value2 = 13.237774030354
$ws.Range("L17").Value2 = 20.304568527918
$ws.Range("M17").Value2 = 24.582560296846
$ws.Range("N17").Value2 = -51.522079172181

# Row 18
$ws.Range("C18").Value2 = 26
$ws.Range("D18").Value2 = 58
$ws.Range("E18").Value2 = -55.172413793103
$ws.Range("F18").Value2 = 139
$ws.Range("G18").Value2 = 195
$ws.Range("H18").Value2 = -28.717948717948
$ws.Range("I18").Value2 = 2306
$ws.Range("J18").Value2 = 2089
$ws.Range("K18").Value2 = 10.387745332695
$ws.Range("L18").Value2 = -12.585291887793
$ws.Range("M18").Value2 = -27.869878010635
$ws.Range("N18").Value2 = -80.651115959053

# Row 19
$ws.Range("C19").Value2 = 90
$ws.Range("D19").Value2 = 109
$ws.Range("E19").Value2 = -17.431192660550
$ws.Range("F19").Value2 = 389
$ws.Range("G19").Value2 = 488
$ws.Range("H19").Value2 = -20.286885245901
$ws.Range("I19").Value2 = 5826
$ws.Range("J19").Value2 = 4791
$ws.Range("K19").Value2 = 21.603005635566
$ws.Range("L19").Value2 = 28.836797877045
$ws.Range("M19").Value2 = 35.268168098444
$ws.Range("N19").Value2 = -15.010940919037

# Row 20
$ws.Range("C20").Value2 = 22
$ws.Range("D20").Value2 = 32
$ws.Range("E20").Value2 = -31.25
$ws.Range("F20").Value2 = 145
$ws.Range("G20").Value2 = 108
$ws.Range("H20").Value2 = 34.259259259259
$ws.Range("I20").Value2 = 1841
$ws.Range("J20").Value2 = 1510
$ws.Range("K20").Value2 = 21.920529801324
$ws.Range("L20").Value2 = 34.281546316557
$ws.Range("M20").Value2 = 30.106007067137
$ws.Range("N20").Value2 = -80.499947039508

# Row 21
$ws.Range("C21").Value2 = 246
$ws.Range("D21").Value2 = 312
$ws.Range("E21").Value2 = -21.153846153846
$ws.Range("F21").Value2 = 1107
$ws.Range("G21").Value2 = 1281
$ws.Range("H21").Value2 = -13.583138173302
$ws.Range("I21").Value2 = 16823
$ws.Range("J21").Value2 = 14339
$ws.Range("K21").Value2 = 17.323383778506
$ws.Range("L21").Value2 = 17.446244065903
$ws.Range("M21").Value2 = 3.986895784398
$ws.Range("N21").Value2 = -68.962972529195

# Row 22
$ws.Range("C22").Value2 = 9
$ws.Range("D22").Value2 = 8
$ws.Range("E22").Value2 = 12.5
$ws.Range("G22").Value2 = 37
$ws.Range("H22").Value2 = -27.027027027027
$ws.Range("I22").Value2 = 344
$ws.Range("J22").Value2 = 289
$ws.Range("K22").Value2 = 19.031141868512
$ws.Range("L22").Value2 = 11.326860841423
$ws.Range("M22").Value2 = -22.171945701357

# Row 23
$ws.Range("C23").Value2 = 26
$ws.Range("D23").Value2 = 24
$ws.Range("E23").Value2 = 8.333333333333
$ws.Range("F23").Value2 = 100
$ws.Range("H23").Value2 = -15.254237288135
$ws.Range("I23").Value2 = 1477
$ws.Range("J23").Value2 = 1455
$ws.Range("K23").Value2 = 1.512027491408
$ws.Range("L23").Value2 = 7.967836257309
$ws.Range("M23").Value2 = 26.999140154772

# Row 24
$ws.Range("C24").Value2 = 188
$ws.Range("D24").Value2 = 199
$ws.Range("E24").Value2 = -5.527638190954
$ws.Range("F24").Value2 = 961
$ws.Range("G24").Value2 = 946
$ws.Range("H24").Value2 = 1.585623678646
$ws.Range("I24").Value2 = 13184
$ws.Range("J24").Value2 = 10479
$ws.Range("K24").Value2 = 25.813531825555
$ws.Range("L24").Value2 = 24.777588491387
$ws.Range("M24").Value2 = 27.001252287833

# Row 25
$ws.Range("C25").Value2 = 84
$ws.Range("D25").Value2 = 83
$ws.Range("E25").Value2 = 1.204819277108
$ws.Range("F25").Value2 = 375
$ws.Range("G25").Value2 = 417
$ws.Range("H25").Value2 = -10.071942446043
$ws.Range("I25").Value2 = 5786
$ws.Range("J25").Value2 = 4780
$ws.Range("K25").Value2 = 21.046025104602
$ws.Range("L25").Value2 = 28.720800889877
$ws.Range("M25").Value2 = -25.090626618332

# Row 26
$ws.Range("C26").Value2 = 9
$ws.Range("D26").Value2 = 4
$ws.Range("E26").Value2 = 125
$ws.Range("F26").Value2 = 26
$ws.Range("G26").Value2 = 23
$ws.Range("H26").Value2 = 13.043478260869
$ws.Range("I26").Value2 = 376
$ws.Range("J26").Value2 = 371
$ws.Range("K26").Value2 = 1.347708894878
$ws.Range("L26").Value2 = 8.985507246376

# Row 27
$ws.Range("C27").Value2 = 11
$ws.Range("D27").Value2 = 11
$ws.Range("E27").Value2 = 0
$ws.Range("F27").Value2 = 37
$ws.Range("G27").Value2 = 62
$ws.Range("H27").Value2 = -40.322580645161
$ws.Range("I27").Value2 = 603
$ws.Range("J27").Value2 = 688
$ws.Range("K27").Value2 = -12.354651162790
$ws.Range("L27").Value2 = 3.965517241379

# Row 28
$ws.Range("C28").Value2 = 6
$ws.Range("D28").Value2 = 2
$ws.Range("E28").Value2 = 200
$ws.Range("F28").Value2 = 20
$ws.Range("G28").Value2 = 22
$ws.Range("H28").Value2 = -9.090909090909
$ws.Range("I28").Value2 = 338
$ws.Range("J28").Value2 = 405
$ws.Range("K28").Value2 = -16.543209876543
$ws.Range("L28").Value2 = -32.128514056224
$ws.Range("M28").Value2 = -31.717171717171
$ws.Range("N28").Value2 = -81.550218340611

# Row 29
$ws.Range("C29").Value2 = 6
$ws.Range("D29").Value2 = 2
$ws.Range("E29").Value2 = 200
$ws.Range("F29").Value2 = 19
$ws.Range("G29").Value2 = 21
$ws.Range("H29").Value2 = -9.523809523809
$ws.Range("I29").Value2 = 286
$ws.Range("J29").Value2 = 333
$ws.Range("K29").Value2 = -14.114114114114
$ws.Range("L29").Value2 = -30.750605326876
$ws.Range("M29").Value2 = -29.032258064516
$ws.Range("N29").Value2 = -82.666666666666

# Row 30
$ws.Range("G30").Value2 = 4
$ws.Range("H30").Value2 = -25
$ws.Range("J30").Value2 = 60
$ws.Range("K30").Value2 = 45

Write-Output "Applied weekly CompStat refresh: header updated to Volume 29 Number 51, week of 12/19/2022-12/25/2022, and precinct/table stats for rows 14-30."
